$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 ("Inicio de sesión"): move status mark from D5 (Parcialmente terminado)
# to C5 (Terminado) -> login feature is now finished.
$ws.Range("D5").Value = $null
$ws.Range("C5").Value = "x"

# Row 6 ("Registro de usuario"): move status mark from C6 (Terminado)
# to D6 (Parcialmente terminado) -> registration is only partially finished now.
$ws.Range("C6").Value = $null
$ws.Range("D6").Value = "x"

# Reflect the new active selection cell left by the author after editing.
$ws.Range("C6").Select()
